# DeveloperGuide: update section of UndoRedoStack to UndoRedoCareTaker
#
# The UndoRedoStack (an "UndoRedo Stack" rectangle fed by a "1"-labelled
# connector arrow) is being removed from the Logic-component class
# diagram, since UndoRedoStack was replaced by UndoRedoCareTaker and
# moved to the Model component.
#
# On the slide, this is represented by three shapes that form one
# visual unit:
#   - "Rectangle 62"                  (id 59) -> the "UndoRedo / Stack" box
#   - "Straight Arrow Connector 57"   (id 61) -> arrow pointing into it
#   - "TextBox 62"                    (id 63) -> the "1" multiplicity label
#
# Remove all three so the diagram no longer references UndoRedoStack.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$idsToRemove = @(59, 61, 63)

for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
    $shp = $s.Shapes.Item($i)
    if ($idsToRemove -contains $shp.Id) {
        $shp.Delete()
    }
}
